$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 updates
$ws.Range("B11").Value = 5
$ws.Range("C11").Value = -1.2

# Row 12 updates
$ws.Range("B12").Value = 105
$ws.Range("C12").Value = -6
$ws.Range("E12").Value = "99.0/140"
